$d = $word.ActiveDocument

# Fix the typo in the title: "FREQUENCY FILTRING" -> "FREQUENCY FILTERING"
$d.Content.Find.Execute("FILTRING", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FILTERING", 2)
